$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing sample sizes (n) for rows 111-112 ---
$ws.Range("I111").Value = 1163
$ws.Range("I112").Value = 1192

# --- New row 114: opinionway poll (partially included in unsure handling) ---
$ws.Range("A114").Value = 33
$ws.Range("B114").Value = 2021
$ws.Range("C114").Value = 12
$ws.Range("D114").Value = 11
$ws.Range("E114").Value = 16
$ws.Range("F114").Value = "opinionway"
$ws.Range("G114").Value = "online"
$ws.Range("H114").Value = "included"
$ws.Range("I114").Value = 1163
$ws.Range("J114").Value = 2
$ws.Range("K114").Value = 1
$ws.Range("L114").Value = 9
$ws.Range("M114").Value = 2
$ws.Range("N114").Value = 3
$ws.Range("O114").Value = 8
$ws.Range("P114").Value = 6
$ws.Range("Q114").Value = 25
$ws.Range("T114").Value = 14
$ws.Range("V114").Value = 4
$ws.Range("W114").Value = 26

# --- New rows 115-120: bva poll (11/19) ---
$ws.Range("A115").Value = 34
$ws.Range("B115").Value = 2021
$ws.Range("C115").Value = 12
$ws.Range("D115").Value = 11
$ws.Range("E115").Value = 14
$ws.Range("F115").Value = "bva"
$ws.Range("G115").Value = "online"
$ws.Range("H115").Value = "excluded"
$ws.Range("I115").Value = 886
$ws.Range("J115").Value = 1.5
$ws.Range("K115").Value = 1
$ws.Range("L115").Value = 7.5
$ws.Range("M115").Value = 2.5
$ws.Range("N115").Value = 2
$ws.Range("O115").Value = 8
$ws.Range("P115").Value = 5
$ws.Range("Q115").Value = 26
$ws.Range("R115").Value = 9
$ws.Range("U115").Value = 1
$ws.Range("V115").Value = 2.5
$ws.Range("W115").Value = 18
$ws.Range("X115").Value = 15
$ws.Range("Y115").Value = 1

$ws.Range("A116").Value = 34
$ws.Range("B116").Value = 2021
$ws.Range("C116").Value = 12
$ws.Range("D116").Value = 11
$ws.Range("E116").Value = 14
$ws.Range("F116").Value = "bva"
$ws.Range("G116").Value = "online"
$ws.Range("H116").Value = "excluded"
$ws.Range("I116").Value = 877
$ws.Range("J116").Value = 1.5
$ws.Range("K116").Value = 1
$ws.Range("L116").Value = 7.5
$ws.Range("M116").Value = 2.5
$ws.Range("N116").Value = 2
$ws.Range("O116").Value = 8.5
$ws.Range("P116").Value = 4
$ws.Range("Q116").Value = 27
$ws.Range("S116").Value = 8
$ws.Range("U116").Value = 1
$ws.Range("V116").Value = 3
$ws.Range("W116").Value = 19
$ws.Range("X116").Value = 15
$ws.Range("Y116").Value = "T_0.5"

$ws.Range("A117").Value = 34
$ws.Range("B117").Value = 2021
$ws.Range("C117").Value = 12
$ws.Range("D117").Value = 11
$ws.Range("E117").Value = 14
$ws.Range("F117").Value = "bva"
$ws.Range("G117").Value = "online"
$ws.Range("H117").Value = "excluded"
$ws.Range("I117").Value = 888
$ws.Range("J117").Value = 1.5
$ws.Range("K117").Value = 1.5
$ws.Range("L117").Value = 7
$ws.Range("M117").Value = 2
$ws.Range("N117").Value = 2
$ws.Range("O117").Value = 8
$ws.Range("P117").Value = 4
$ws.Range("Q117").Value = 26
$ws.Range("T117").Value = 12
$ws.Range("U117").Value = 1
$ws.Range("V117").Value = 1.5
$ws.Range("W117").Value = 18
$ws.Range("X117").Value = 15
$ws.Range("Y117").Value = 0.5

$ws.Range("A118").Value = 34
$ws.Range("B118").Value = 2021
$ws.Range("C118").Value = 12
$ws.Range("D118").Value = 11
$ws.Range("E118").Value = 14
$ws.Range("F118").Value = "bva"
$ws.Range("G118").Value = "online"
$ws.Range("H118").Value = "excluded"
$ws.Range("I118").Value = 870
$ws.Range("J118").Value = 1.5
$ws.Range("K118").Value = 1.5
$ws.Range("L118").Value = 8
$ws.Range("M118").Value = 2
$ws.Range("N118").Value = 2
$ws.Range("O118").Value = 9
$ws.Range("P118").Value = 4
$ws.Range("Q118").Value = 28
$ws.Range("R118").Value = 11
$ws.Range("U118").Value = 1.5
$ws.Range("V118").Value = 4
$ws.Range("W118").Value = 27
$ws.Range("Y118").Value = 0.5

$ws.Range("A119").Value = 34
$ws.Range("B119").Value = 2021
$ws.Range("C119").Value = 12
$ws.Range("D119").Value = 11
$ws.Range("E119").Value = 14
$ws.Range("F119").Value = "bva"
$ws.Range("G119").Value = "online"
$ws.Range("H119").Value = "excluded"
$ws.Range("I119").Value = 856
$ws.Range("J119").Value = 1.5
$ws.Range("K119").Value = 1.5
$ws.Range("L119").Value = 8
$ws.Range("M119").Value = 2
$ws.Range("N119").Value = 2
$ws.Range("O119").Value = 8
$ws.Range("P119").Value = 5
$ws.Range("Q119").Value = 29
$ws.Range("S119").Value = 9
$ws.Range("U119").Value = 1.5
$ws.Range("V119").Value = 4
$ws.Range("W119").Value = 28
$ws.Range("Y119").Value = 0.5

$ws.Range("A120").Value = 34
$ws.Range("B120").Value = 2021
$ws.Range("C120").Value = 12
$ws.Range("D120").Value = 11
$ws.Range("E120").Value = 14
$ws.Range("F120").Value = "bva"
$ws.Range("G120").Value = "online"
$ws.Range("H120").Value = "excluded"
$ws.Range("I120").Value = 873
$ws.Range("J120").Value = 1
$ws.Range("K120").Value = 1
$ws.Range("L120").Value = 7.5
$ws.Range("M120").Value = 2
$ws.Range("N120").Value = 2
$ws.Range("O120").Value = 8.5
$ws.Range("P120").Value = 5
$ws.Range("Q120").Value = 27
$ws.Range("T120").Value = 15
$ws.Range("U120").Value = 1.5
$ws.Range("V120").Value = 3
$ws.Range("W120").Value = 26
$ws.Range("Y120").Value = 0.5

# --- Restore view state (best effort) ---
$ws.Range("Y114").Select()
